$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.249.36"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "1.673.19"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.67"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5136"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2665"
$ws.Range("E8").Value = "  +3.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06388"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.60"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07386"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "1.674.07"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.557"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5841"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "1.900.17"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008679"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.70"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "26.321.96"
$ws.Range("E18").Value = "  -2.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.972"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.90"
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "189.79"
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.222"
$ws.Range("E23").Value = "  -1.12%  "
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.09"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.670"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1184"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.69"
$ws.Range("E28").Value = "  +3.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06013"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.284"
$ws.Range("E30").Value = "  -3.95%  "
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.534"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.535"
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.647"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.017"
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6041"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.374"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.091"
$ws.Range("E40").Value = "  +4.56%  "
$ws.Range("D41").Value = "1.083.65"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8691"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.35"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").Value = "1.820.21"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000112"
$ws.Range("E46").Value = "  +8.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.46"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.012"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.055"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05215"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4296"
$ws.Range("E51").Value = "  -1.76%  "
